{"js": "// Replace the division-problem text in each table cell with the new\n// problem, matching the document's original left-to-right, top-to-bottom\n// order. Each \"from\" string is unique in the document, so a plain\n// case-sensitive search locates the single correct cell. All searches are\n// issued (and loaded) against the pristine document BEFORE any text is\n// replaced, so earlier replacements that happen to produce text equal to\n// a later \"from\" value (e.g. \"57\u00f72=\") cannot be re-matched by mistake.\nconst replacements = [\n  [\"88\u00f78=\", \"80\u00f72=\"],\n  [\"76\u00f79=\", \"28\u00f76=\"],\n  [\"73\u00f77=\", \"48\u00f76=\"],\n  [\"50\u00f76=\", \"57\u00f72=\"],\n  [\"67\u00f77=\", \"16\u00f77=\"],\n  [\"96\u00f77=\", \"55\u00f75=\"],\n  [\"46\u00f77=\", \"33\u00f74=\"],\n  [\"63\u00f74=\", \"91\u00f73=\"],\n  [\"96\u00f76=\", \"98\u00f72=\"],\n  [\"19\u00f73=\", \"64\u00f73=\"],\n  [\"92\u00f74=\", \"16\u00f77=\"],\n  [\"81\u00f79=\", \"20\u00f74=\"],\n  [\"90\u00f75=\", \"96\u00f73=\"],\n  [\"91\u00f73=\", \"71\u00f73=\"],\n  [\"11\u00f72=\", \"89\u00f73=\"],\n  [\"14\u00f73=\", \"48\u00f76=\"],\n  [\"59\u00f76=\", \"51\u00f74=\"],\n  [\"38\u00f77=\", \"23\u00f76=\"],\n  [\"56\u00f78=\", \"63\u00f77=\"],\n  [\"93\u00f75=\", \"45\u00f79=\"],\n  [\"88\u00f73=\", \"12\u00f74=\"],\n  [\"57\u00f72=\", \"24\u00f73=\"],\n  [\"15\u00f79=\", \"36\u00f72=\"],\n  [\"37\u00f73=\", \"70\u00f75=\"],\n  [\"57\u00f74=\", \"85\u00f72=\"],\n];\n\nconst body = context.document.body;\nconst searchResults = [];\nfor (const [from] of replacements) {\n  const results = body.search(from, { matchCase: true, matchWholeWord: false });\n  results.load(\"items/text\");\n  searchResults.push(results);\n}\nawait context.sync();\n\nfor (let i = 0; i < replacements.length; i++) {\n  const [from, to] = replacements[i];\n  const results = searchResults[i];\n  for (const item of results.items) {\n    item.insertText(to, Word.InsertLocation.replace);\n  }\n}\nawait context.sync();\n", "ps1": "# Replace the division-problem text in each table cell with the new\n# problem, matching the document's original left-to-right, top-to-bottom\n# order. Each \"from\" string is unique in the document, so Find.Execute\n# locates exactly the right cell.\n#\n# All target cells are located (via Find, on a duplicated Range so the\n# main document Range is untouched) BEFORE any text is written back, and\n# only then are the captured Range.Text properties set. That two-phase\n# approach is required because some \"to\" values equal a later \"from\"\n# value (e.g. \"57\u00f72=\" is written at step 4 and is itself the thing being\n# searched for at step 22) \u2014 doing find-then-replace-immediately in a\n# single pass could re-match a cell that was only just rewritten.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{from = \"88\u00f78=\"; to = \"80\u00f72=\"},\n    @{from = \"76\u00f79=\"; to = \"28\u00f76=\"},\n    @{from = \"73\u00f77=\"; to = \"48\u00f76=\"},\n    @{from = \"50\u00f76=\"; to = \"57\u00f72=\"},\n    @{from = \"67\u00f77=\"; to = \"16\u00f77=\"},\n    @{from = \"96\u00f77=\"; to = \"55\u00f75=\"},\n    @{from = \"46\u00f77=\"; to = \"33\u00f74=\"},\n    @{from = \"63\u00f74=\"; to = \"91\u00f73=\"},\n    @{from = \"96\u00f76=\"; to = \"98\u00f72=\"},\n    @{from = \"19\u00f73=\"; to = \"64\u00f73=\"},\n    @{from = \"92\u00f74=\"; to = \"16\u00f77=\"},\n    @{from = \"81\u00f79=\"; to = \"20\u00f74=\"},\n    @{from = \"90\u00f75=\"; to = \"96\u00f73=\"},\n    @{from = \"91\u00f73=\"; to = \"71\u00f73=\"},\n    @{from = \"11\u00f72=\"; to = \"89\u00f73=\"},\n    @{from = \"14\u00f73=\"; to = \"48\u00f76=\"},\n    @{from = \"59\u00f76=\"; to = \"51\u00f74=\"},\n    @{from = \"38\u00f77=\"; to = \"23\u00f76=\"},\n    @{from = \"56\u00f78=\"; to = \"63\u00f77=\"},\n    @{from = \"93\u00f75=\"; to = \"45\u00f79=\"},\n    @{from = \"88\u00f73=\"; to = \"12\u00f74=\"},\n    @{from = \"57\u00f72=\"; to = \"24\u00f73=\"},\n    @{from = \"15\u00f79=\"; to = \"36\u00f72=\"},\n    @{from = \"37\u00f73=\"; to = \"70\u00f75=\"},\n    @{from = \"57\u00f74=\"; to = \"85\u00f72=\"}\n)\n\n$targetRanges = @()\nforeach ($item in $replacements) {\n    $searchRange = $d.Content.Duplicate\n    $searchRange.Find.Execute($item.from, $false, $false, $false, $false, $false, $true, 1, $false, \"\", 0) | Out-Null\n    $targetRanges += ,$searchRange.Duplicate\n}\n\nfor ($i = 0; $i -lt $replacements.Count; $i++) {\n    $targetRanges[$i].Text = $replacements[$i].to\n}\n"}
